$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date serial that needs to move from
# 2023-09-21 (serial 45190) to 2023-09-23 (serial 45192) for every data
# row (rows 2 through 123).
$ws.Range("C2:C123").Value2 = 45192
